$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 374
$ws.Range("I2").Value = 944
$ws.Range("J2").Value = 3657
$ws.Range("K2").Value = 15
$ws.Range("L2").Value = 1021
$ws.Range("M2").Value = 53
$ws.Range("N2").Value = 679
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 14
$ws.Range("Q2").Value = 9
$ws.Range("R2").Value = 41
$ws.Range("S2").Value = 384
$ws.Range("T2").Value = 683
$ws.Range("U2").Value = 45
$ws.Range("V2").Value = 5983
$ws.Range("X2").Value = 5896
$ws.Range("Y2").Value = 13
$ws.Range("Z2").Value = 90
$ws.Range("AA2").Value = 27
